$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "# records" count for the "C, Si, N" section (row 41), mirroring
#     the other section headers (e.g. row 33 / row 48) that carry a record
#     count in column C and the "# records" label in column D.
$ws.Cells.Item(41, 3).Value = 1301
$ws.Cells.Item(41, 3).NumberFormat = "#,##0"
$ws.Cells.Item(41, 4).Value = "# records"

# --- Insert a new row for the XGBoost model result inside the "C, Si, N"
#     section, between kNN (95.7) and Random Forest Classifier (96.9), since
#     that section's rows are kept sorted ascending by rating.
$ws.Rows.Item(46).Insert()
$ws.Cells.Item(46, 1).Value = "XGBoost"
$ws.Cells.Item(46, 2).Value = 96.3

# --- The section's sort range grew by the inserted row; extend it to match.
$sortObj = $ws.Sort
$sortObj.SortFields.Add($ws.Range("B42:B47")) | Out-Null
$sortObj.SetRange($ws.Range("A42:B47")) | Out-Null
$sortObj.Apply() | Out-Null

# --- Update the active selection to reflect where the edit was made.
$ws.Range("D41").Select() | Out-Null
